# Add data for 2022-12-19: update the "as of" date labels and the
# December 2022 / Total figures in the "Total" column (I).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date.
$ws.Name = "Through 2022-12-11"

# Update the header label in I1 ("2022 (through 12-10)" -> "... 12-11").
$ws.Range("I1").Value = "2022 (through 12-11)"

# December 2022 total (row 13) goes from 44 to 50.
$ws.Range("I13").Value = 50

# Grand total (row 14) goes from 1560 to 1566.
$ws.Range("I14").Value = 1566
